$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03208554907185
$ws.Range("D2").Value = 1.036464701325385
$ws.Range("E2").Value = 1.031565745319425
$ws.Range("F2").Value = 1.040925758230683
$ws.Range("I2").Value = 1.036856252633784
$ws.Range("J2").Value = 1.037217709288813
$ws.Range("K2").Value = 1.039258366917583
$ws.Range("L2").Value = 1.034373494534126
$ws.Range("M2").Value = 1.043706724281097
$ws.Range("N2").Value = 1.016343574211428

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033060605156285
$ws.Range("D3").Value = 1.037225263268745
$ws.Range("E3").Value = 1.032395492424884
$ws.Range("F3").Value = 1.043154079527121
$ws.Range("I3").Value = 1.037163262644874
$ws.Range("J3").Value = 1.037834782222525
$ws.Range("K3").Value = 1.039828665962222
$ws.Range("L3").Value = 1.03501177531797
$ws.Range("M3").Value = 1.04574184721182
$ws.Range("N3").Value = 1.016552437025543

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033691436801305
$ws.Range("D4").Value = 1.037717193137872
$ws.Range("E4").Value = 1.032932661753921
$ws.Range("F4").Value = 1.04459057057502
$ws.Range("I4").Value = 1.037360385180794
$ws.Range("J4").Value = 1.038233380728532
$ws.Range("K4").Value = 1.040196822972799
$ws.Range("L4").Value = 1.035424408265053
$ws.Range("M4").Value = 1.047053002973725
$ws.Range("N4").Value = 1.016687252024688

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033956616738723
$ws.Range("D5").Value = 1.037923952154809
$ws.Range("E5").Value = 1.033158552096423
$ws.Range("F5").Value = 1.045193215283441
$ws.Range("I5").Value = 1.037442889893925
$ws.Range("J5").Value = 1.038400787423748
$ws.Range("K5").Value = 1.040351390136164
$ws.Range("L5").Value = 1.035597789361537
$ws.Range("M5").Value = 1.047602877182719
$ws.Range("N5").Value = 1.016743848660431

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034001140374573
$ws.Range("D6").Value = 1.037958665090237
$ws.Range("E6").Value = 1.033196483861783
$ws.Range("F6").Value = 1.045294329252063
$ws.Range("I6").Value = 1.037456721402795
$ws.Range("J6").Value = 1.038428886143951
$ws.Range("K6").Value = 1.040377330587475
$ws.Range("L6").Value = 1.035626895560714
$ws.Range("M6").Value = 1.047695126023908
$ws.Range("N6").Value = 1.016753346822515

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033694980232587
$ws.Range("D7").Value = 1.037719956050979
$ws.Range("E7").Value = 1.032935679858498
$ws.Range("F7").Value = 1.04459862803955
$ws.Range("I7").Value = 1.037361489047536
$ws.Range("J7").Value = 1.038235618268879
$ws.Range("K7").Value = 1.040198889115286
$ws.Range("L7").Value = 1.035426725343931
$ws.Range("M7").Value = 1.047060355631539
$ws.Range("N7").Value = 1.016688008584229

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032415094098977
$ws.Range("D8").Value = 1.036721779565704
$ws.Range("E8").Value = 1.031846107059873
$ws.Range("F8").Value = 1.04167996375131
$ws.Range("I8").Value = 1.036960326539478
$ws.Range("J8").Value = 1.037426395125926
$ws.Range("K8").Value = 1.039451281748593
$ws.Range("L8").Value = 1.034589282998067
$ws.Range("M8").Value = 1.044395703054612
$ws.Range("N8").Value = 1.016414229481609

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03015899758677
$ws.Range("D9").Value = 1.034961266326452
$ws.Range("E9").Value = 1.029928172868494
$ws.Range("F9").Value = 1.036494223633649
$ws.Range("I9").Value = 1.036241613162507
$ws.Range("J9").Value = 1.035995125924091
$ws.Range("K9").Value = 1.038127236314369
$ws.Range("L9").Value = 1.033110685784543
$ws.Range("M9").Value = 1.03965519003275
$ws.Range("N9").Value = 1.015929235084259

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028654342798618
$ws.Range("D10").Value = 1.033786466836914
$ws.Range("E10").Value = 1.028650889568175
$ws.Range("F10").Value = 1.033006306175898
$ws.Range("I10").Value = 1.035754428627578
$ws.Range("J10").Value = 1.035037316330466
$ws.Range("K10").Value = 1.03723999279155
$ws.Range("L10").Value = 1.032122952792801
$ws.Range("M10").Value = 1.036462635075604
$ws.Range("N10").Value = 1.015604170377068

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028002653662186
$ws.Range("D11").Value = 1.033277486839218
$ws.Range("E11").Value = 1.028098122217212
$ws.Range("F11").Value = 1.031488240686428
$ws.Range("I11").Value = 1.035541541923221
$ws.Range("J11").Value = 1.03462169944087
$ws.Range("K11").Value = 1.036854712634316
$ws.Range("L11").Value = 1.031694768680727
$ws.Range("M11").Value = 1.035072151067889
$ws.Range("N11").Value = 1.015462998804442

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027760561118529
$ws.Range("D12").Value = 1.033088385147969
$ws.Range("E12").Value = 1.027892844847037
$ws.Range("F12").Value = 1.030923155361263
$ws.Range("I12").Value = 1.035462173916887
$ws.Range("J12").Value = 1.034467187467568
$ws.Range("K12").Value = 1.036711436041775
$ws.Range("L12").Value = 1.031535647562483
$ws.Range("M12").Value = 1.034554410963291
$ws.Range("N12").Value = 1.015410498498324

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027812492001655
$ws.Range("D13").Value = 1.033128950106709
$ws.Range("E13").Value = 1.027936875514186
$ws.Range("F13").Value = 1.031044423408784
$ws.Range("I13").Value = 1.03547921188909
$ws.Range("J13").Value = 1.034500336866084
$ws.Range("K13").Value = 1.036742176910268
$ws.Range("L13").Value = 1.03156978296862
$ws.Range("M13").Value = 1.034665525257342
$ws.Range("N13").Value = 1.015421762844038

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027982642760441
$ws.Range("D14").Value = 1.033261856536093
$ws.Range("E14").Value = 1.028081153011494
$ws.Range("F14").Value = 1.03144155543947
$ws.Range("I14").Value = 1.03553498731609
$ws.Range("J14").Value = 1.034608930160971
$ws.Range("K14").Value = 1.036842872749249
$ws.Range("L14").Value = 1.031681617204468
$ws.Range("M14").Value = 1.035029380277876
$ws.Range("N14").Value = 1.015458660397867

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02808747474105
$ws.Range("D15").Value = 1.033343738726093
$ws.Range("E15").Value = 1.028170053084043
$ws.Range("F15").Value = 1.031686080275814
$ws.Range("I15").Value = 1.035569313588672
$ws.Range("J15").Value = 1.03467582036383
$ws.Range("K15").Value = 1.036904892687102
$ws.Range("L15").Value = 1.031750512069866
$ws.Range("M15").Value = 1.035253396272146
$ws.Range("N15").Value = 1.015481385853594

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02869758967678
$ws.Range("D16").Value = 1.033820240078364
$ws.Range("E16").Value = 1.028687581274301
$ws.Range("F16").Value = 1.033106887971392
$ws.Range("I16").Value = 1.035768516325084
$ws.Range("J16").Value = 1.035064880825229
$ws.Range("K16").Value = 1.037265539313494
$ws.Range("L16").Value = 1.032151359572494
$ws.Range("M16").Value = 1.036554743356663
$ws.Range("N16").Value = 1.015613530667608

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029080253722361
$ws.Range("D17").Value = 1.034119059914762
$ws.Range("E17").Value = 1.029012294106269
$ws.Range("F17").Value = 1.033996013542644
$ws.Range("I17").Value = 1.035892952091841
$ws.Range("J17").Value = 1.035308691934459
$ws.Range("K17").Value = 1.037491468476981
$ws.Range("L17").Value = 1.032402668966351
$ws.Range("M17").Value = 1.03736885270019
$ws.Range("N17").Value = 1.015696309899057

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029303439503425
$ws.Range("D18").Value = 1.034293329144429
$ws.Range("E18").Value = 1.029201723064678
$ws.Range("F18").Value = 1.034513878097117
$ws.Range("I18").Value = 1.035965347044289
$ws.Range("J18").Value = 1.035450818107408
$ws.Range("K18").Value = 1.03762314330895
$ws.Range("L18").Value = 1.032549206350993
$ws.Range("M18").Value = 1.037842931227608
$ws.Range("N18").Value = 1.015744553470585

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029379537474183
$ws.Range("D19").Value = 1.034352745848045
$ws.Range("E19").Value = 1.029266318462364
$ws.Range("F19").Value = 1.034690330971317
$ws.Range("I19").Value = 1.035990000309689
$ws.Range("J19").Value = 1.035499265127425
$ws.Range("K19").Value = 1.037668023116706
$ws.Range("L19").Value = 1.03259916384641
$ws.Range("M19").Value = 1.03800444897193
$ws.Range("N19").Value = 1.01576099648089

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029039199125516
$ws.Range("D20").Value = 1.034087002209295
$ws.Range("E20").Value = 1.028977452451061
$ws.Range("F20").Value = 1.033900696374348
$ws.Range("I20").Value = 1.03587962059903
$ws.Range("J20").Value = 1.035282542082456
$ws.Range("K20").Value = 1.03746723937256
$ws.Range("L20").Value = 1.032375710729067
$ws.Range("M20").Value = 1.037281587180891
$ws.Range("N20").Value = 1.015687432626682

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027932538338578
$ws.Range("D21").Value = 1.03322272012562
$ws.Range("E21").Value = 1.028038665647532
$ws.Range("F21").Value = 1.031324643633627
$ws.Range("I21").Value = 1.03551857093155
$ws.Range("J21").Value = 1.03457695583406
$ws.Range("K21").Value = 1.036813224939371
$ws.Range("L21").Value = 1.031648686873256
$ws.Range("M21").Value = 1.034922268920827
$ws.Range("N21").Value = 1.015447796726987

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027236584483221
$ws.Range("D22").Value = 1.032679057313196
$ws.Range("E22").Value = 1.027448673291266
$ws.Range("F22").Value = 1.029697966273709
$ws.Range("I22").Value = 1.035289871857939
$ws.Range("J22").Value = 1.034132553536413
$ws.Range("K22").Value = 1.036401056694964
$ws.Range("L22").Value = 1.031191146727816
$ws.Range("M22").Value = 1.033431606652891
$ws.Range("N22").Value = 1.015296764038695

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027605537765608
$ws.Range("D23").Value = 1.032967287771751
$ws.Range("E23").Value = 1.027761415000084
$ws.Range("F23").Value = 1.030560976867671
$ws.Range("I23").Value = 1.035411270720841
$ws.Range("J23").Value = 1.034368213240825
$ws.Range("K23").Value = 1.036619646693121
$ws.Range("L23").Value = 1.031433738640806
$ws.Range("M23").Value = 1.034222536275855
$ws.Range("N23").Value = 1.015376863951459

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029057749979557
$ws.Range("D24").Value = 1.034101487791277
$ws.Range("E24").Value = 1.028993195804695
$ws.Range("F24").Value = 1.033943768407536
$ws.Range("I24").Value = 1.035885645102604
$ws.Range("J24").Value = 1.035294358337041
$ws.Range("K24").Value = 1.037478187789171
$ws.Range("L24").Value = 1.032387892142039
$ws.Range("M24").Value = 1.037321021114307
$ws.Range("N24").Value = 1.015691444008132

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030742351415773
$ws.Range("D25").Value = 1.035416595577486
$ws.Range("E25").Value = 1.030423766651335
$ws.Range("F25").Value = 1.037840129836832
$ws.Range("I25").Value = 1.036428827894004
$ws.Range("J25").Value = 1.036365778874097
$ws.Range("K25").Value = 1.03847033027991
$ws.Range("L25").Value = 1.033493288320661
$ws.Range("M25").Value = 1.040886261345198
$ws.Range("N25").Value = 1.016054922558215
